$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '43.062.37'
$ws.Range('E2').Value = '  +0.88%  '
$ws.Range('D3').Value = '2.282.43'
$ws.Range('E3').Value = '  +1.80%  '
$ws.Range('E4').Value = '  +0.39%  '
$ws.Range('D5').NumberFormat = "@"
$ws.Range('D5').Value = '113.10'
$ws.Range('E5').Value = '  -1.96%  '
$ws.Range('D6').NumberFormat = "@"
$ws.Range('D6').Value = '308.16'
$ws.Range('E6').Value = '  +6.75%  '
$ws.Range('D7').NumberFormat = "@"
$ws.Range('D7').Value = '0.631'
$ws.Range('E7').Value = '  +0.36%  '
$ws.Range('E8').Value = '  -0.14%  '
$ws.Range('D9').NumberFormat = "@"
$ws.Range('D9').Value = '0.613'
$ws.Range('E9').Value = '  -0.54%  '
$ws.Range('D10').NumberFormat = "@"
$ws.Range('D10').Value = '44.57'
$ws.Range('E10').Value = '  -4.78%  '
$ws.Range('D11').NumberFormat = "@"
$ws.Range('D11').Value = '0.0926'
$ws.Range('E11').Value = '  -0.86%  '
$ws.Range('D12').NumberFormat = "@"
$ws.Range('D12').Value = '55.12'
$ws.Range('E12').Value = '  +0.84%  '
$ws.Range('D13').NumberFormat = "@"
$ws.Range('D13').Value = '8.85'
$ws.Range('E13').Value = '  -4.01%  '
$ws.Range('E14').Value = '  +18.34%  '
$ws.Range('E15').Value = '  -0.05%  '
$ws.Range('D16').NumberFormat = "@"
$ws.Range('D16').Value = '15.44'
$ws.Range('E16').Value = '  -0.10%  '
$ws.Range('D17').Value = '2.624.48'
$ws.Range('E17').Value = '  +1.65%  '
$ws.Range('D18').Value = '2.282.55'
$ws.Range('E18').Value = '  +1.46%  '
$ws.Range('D19').Value = '43.025.25'
$ws.Range('E19').Value = '  +0.63%  '
$ws.Range('E20').Value = '  -0.15%  '
$ws.Range('D21').NumberFormat = "@"
$ws.Range('D21').Value = '7.24'
$ws.Range('E21').Value = '  +0.08%  '
$ws.Range('D22').NumberFormat = "@"
$ws.Range('D22').Value = '75.30'
$ws.Range('E22').Value = '  +2.51%  '
$ws.Range('D23').NumberFormat = "@"
$ws.Range('D23').Value = '3.60'
$ws.Range('E23').Value = '  +7.68%  '
$ws.Range('D24').NumberFormat = "@"
$ws.Range('D24').Value = '2.47'
$ws.Range('E24').Value = '  +4.04%  '
$ws.Range('D25').NumberFormat = "@"
$ws.Range('D25').Value = '254.50'
$ws.Range('E25').Value = '  +9.49%  '
$ws.Range('D26').NumberFormat = "@"
$ws.Range('D26').Value = '9.00'
$ws.Range('E26').Value = '  -2.43%  '
$ws.Range('D27').NumberFormat = "@"
$ws.Range('D27').Value = '11.74'
$ws.Range('E27').Value = '  -3.55%  '
$ws.Range('D28').NumberFormat = "@"
$ws.Range('D28').Value = '0.999'
$ws.Range('E28').Value = '  -0.24%  '
$ws.Range('E29').Value = '  +2.39%  '
$ws.Range('D30').NumberFormat = "@"
$ws.Range('D30').Value = '38.30'
$ws.Range('E30').Value = '  -4.87%  '
$ws.Range('B31').Value = 'Monero'
$ws.Range('C31').Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range('D31').NumberFormat = "@"
$ws.Range('D31').Value = '174.87'
$ws.Range('E31').Value = '  -0.26%  '
$ws.Range('B32').Value = 'EthereumClassic'
$ws.Range('C32').Value = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
$ws.Range('D32').NumberFormat = "@"
$ws.Range('D32').Value = '22.17'
$ws.Range('E32').Value = '  +4.15%  '
$ws.Range('D33').NumberFormat = "@"
$ws.Range('D33').Value = '3.18'
$ws.Range('E33').Value = '  -3.45%  '
$ws.Range('D34').NumberFormat = "@"
$ws.Range('D34').Value = '0.0901'
$ws.Range('E34').Value = '  -1.23%  '
$ws.Range('E35').Value = '  +1.82%  '
$ws.Range('D36').NumberFormat = "@"
$ws.Range('D36').Value = '5.07'
$ws.Range('E36').Value = '  +8.99%  '
$ws.Range('E37').Value = '  +0.84%  '
$ws.Range('D38').NumberFormat = "@"
$ws.Range('D38').Value = '4.22'
$ws.Range('E38').Value = '  -7.05%  '
$ws.Range('D39').NumberFormat = "@"
$ws.Range('D39').Value = '0.0378'
$ws.Range('E39').Value = '  +1.02%  '
$ws.Range('E40').Value = '  -1.26%  '
$ws.Range('D41').NumberFormat = "@"
$ws.Range('D41').Value = '2.53'
$ws.Range('E41').Value = '  -4.87%  '
$ws.Range('D42').NumberFormat = "@"
$ws.Range('D42').Value = '72.50'
$ws.Range('E42').Value = '  -0.74%  '
$ws.Range('D43').NumberFormat = "@"
$ws.Range('D43').Value = '0.231'
$ws.Range('E43').Value = '  -2.26%  '
$ws.Range('D45').NumberFormat = "@"
$ws.Range('D45').Value = '12.65'
$ws.Range('E45').Value = '  -6.40%  '
$ws.Range('E46').Value = '  +3.64%  '
$ws.Range('D47').NumberFormat = "@"
$ws.Range('D47').Value = '5.67'
$ws.Range('E47').Value = '  +1.00%  '
$ws.Range('D48').NumberFormat = "@"
$ws.Range('D48').Value = '107.41'
$ws.Range('E48').Value = '  +5.96%  '
$ws.Range('D49').NumberFormat = "@"
$ws.Range('D49').Value = '1.30'
$ws.Range('E49').Value = '  -0.84%  '
$ws.Range('D50').NumberFormat = "@"
$ws.Range('D50').Value = '8.80'
$ws.Range('E50').Value = '  +2.98%  '
$ws.Range('D51').NumberFormat = "@"
$ws.Range('D51').Value = '73.37'
$ws.Range('E51').Value = '  +4.94%  '
